{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2023-10-26 Thursday\", \"2023-10-27 Friday\"],\n  [\"31\u00d719=\", \"11\u00d744=\"],\n  [\"72\u00d720=\", \"81\u00d759=\"],\n  [\"92\u00d792=\", \"87\u00d765=\"],\n  [\"27\u00d760=\", \"66\u00d720=\"],\n  [\"78\u00d727=\", \"25\u00d743=\"],\n  [\"29\u00d726=\", \"80\u00d781=\"],\n  [\"50\u00d716=\", \"61\u00d799=\"],\n  [\"42\u00d761=\", \"67\u00d751=\"],\n  [\"45\u00d795=\", \"55\u00d749=\"],\n  [\"26\u00d779=\", \"97\u00d748=\"],\n  [\"87\u00d799=\", \"55\u00d799=\"],\n  [\"11\u00d727=\", \"63\u00d744=\"],\n  [\"77\u00d758=\", \"85\u00d721=\"],\n  [\"50\u00d789=\", \"63\u00d721=\"],\n  [\"35\u00d754=\", \"88\u00d734=\"],\n  [\"76\u00d798=\", \"51\u00d790=\"],\n  [\"25\u00d769=\", \"23\u00d733=\"],\n  [\"13\u00d723=\", \"45\u00d769=\"],\n  [\"20\u00d735=\", \"40\u00d790=\"],\n  [\"89\u00d718=\", \"37\u00d787=\"],\n  [\"34\u00d727=\", \"32\u00d796=\"],\n  [\"45\u00d755=\", \"64\u00d770=\"],\n  [\"13\u00d735=\", \"21\u00d717=\"],\n  [\"47\u00d749=\", \"75\u00d744=\"],\n  [\"27\u00d724=\", \"55\u00d764=\"],\n];\n\nconst searchResults = [];\nfor (const [oldText, newText] of replacements) {\n  const res = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  res.load(\"items\");\n  searchResults.push([res, newText]);\n}\nawait context.sync();\n\nfor (const [res, newText] of searchResults) {\n  for (let i = 0; i < res.items.length; i++) {\n    res.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2023-10-26 Thursday', '2023-10-27 Friday'),\n    @('31\u00d719=', '11\u00d744='),\n    @('72\u00d720=', '81\u00d759='),\n    @('92\u00d792=', '87\u00d765='),\n    @('27\u00d760=', '66\u00d720='),\n    @('78\u00d727=', '25\u00d743='),\n    @('29\u00d726=', '80\u00d781='),\n    @('50\u00d716=', '61\u00d799='),\n    @('42\u00d761=', '67\u00d751='),\n    @('45\u00d795=', '55\u00d749='),\n    @('26\u00d779=', '97\u00d748='),\n    @('87\u00d799=', '55\u00d799='),\n    @('11\u00d727=', '63\u00d744='),\n    @('77\u00d758=', '85\u00d721='),\n    @('50\u00d789=', '63\u00d721='),\n    @('35\u00d754=', '88\u00d734='),\n    @('76\u00d798=', '51\u00d790='),\n    @('25\u00d769=', '23\u00d733='),\n    @('13\u00d723=', '45\u00d769='),\n    @('20\u00d735=', '40\u00d790='),\n    @('89\u00d718=', '37\u00d787='),\n    @('34\u00d727=', '32\u00d796='),\n    @('45\u00d755=', '64\u00d770='),\n    @('13\u00d735=', '21\u00d717='),\n    @('47\u00d749=', '75\u00d744='),\n    @('27\u00d724=', '55\u00d764='),\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
